$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Vivo', '৳25,000.00', '2024, April 29', 'Android 14, OriginOS 4', '6.78" 1260x2800 pixels', '50MP 2160p', '8/12GB RAM Snapdragon 7 Gen 3', '6000mAh 80W')
    ,('Vivo', '৳35,000.00', '2024, April 29', 'Android 14, OriginOS 4', '6.78" 1260x2800 pixels', '50MP 2160p', '12/16GB RAM Snapdragon 8s Gen 3', '6000mAh 80W')
    ,('Vivo', '৳20,000.00', '2024, April 29', 'Android 14, OriginOS 4', '6.72" 1080x2408 pixels', '50MP 2160p', '8/12GB RAM Snapdragon 6 Gen 1', '6000mAh 44W')
    ,('Realme', '৳18,000.00', '2024, April 24', 'Android 14, Realme UI 5.0', '6.72" 1080x2400 pixels', '50MP 1080p', '6.72" 1080x2400 pixels', '5000mAh 45W')
    ,('Infinix', 'Coming soon', 'Not announced yet', 'Android 13', '6.78" 1080x2400 pixels', '108MP 1440p', '4/8GB RAM Dimensity 7020', '5000mAh Li-Po')
    ,('Realme', '৳20,000.00', '2024, April 25', 'Android 14, Realme UI 5.0', '6.67" 1080x2400 pixels', '50MP 1080p', '6/8GB RAM Dimensity 7050', '5000mAh 45W')
    ,('Oppo', '৳32,000.00', '2024, April 29', 'Android 14, ColorOS 14', '6.7" 1080x2412 pixels', '50MP 2160p', '8/12GB RAM Snapdragon 7 Gen 3', '5500mAh 100W')
    ,('HMD', '৳25,000.00', '2024, April', 'Android 14', '6.56" 720x1480 pixels', '50MP 1080p', '4-8GB RAM Unisoc T606', '5000mAh 20W')
    ,('Nothing', 'Coming soon', 'Not announced yet', 'Android 14', '6.55" 1440x3216 pixels', '50MP 2160p', '8/12GB RAM Snapdragon 778G+ 5G', '4500mAh Li-Po')
    ,('Vivo', 'Coming soon', 'Not announced yet', 'Android 14', '6.78" 1440x3200 pixels', '64MP 4320p', '16GB RAM Snapdragon 8 Gen 3', '5100mAh Li-Po')
    ,('Vivo', 'Coming soon', 'Not announced yet', 'Android 14', '6.78" 1260x2800 pixels', '64MP 4320p', '12/16GB RAM Snapdragon 8 Gen 3', '5000mAh Li-Po')
    ,('Motorola', 'Coming soon', 'Not announced yet', 'Android 13', '6.5" 720x1600 pixels', '13MP 1080p', '2/4GB RAM Unisoc T606', '5000mAh Li-Po')
    ,('Huawei', '৳160,000.00', '2024, April 29', 'HarmonyOS 4.2', '6.8" 1260x2844 pixels', '50MP 2160p', '16GB RAM Kirin 9010', '5200mAh 100W80W')
    ,('Huawei', '৳120,000.00', 'Exp. release 2024, April 29', 'HarmonyOS 4.2', '6.8" 1260x2844 pixels', '50MP 2160p', '16GB RAM Kirin 9010', '5050mAh 100W 80W')
    ,('Huawei', '৳100,000.00', '2024, April 29', 'HarmonyOS 4.2', '6.8" 1260x2844 pixels', '50MP 2160p', '12GB RAM Kirin 9010', '5050mAh 100W 80W')
    ,('Huawei', '৳85,000.00', 'Exp. release 2024, May 06', 'EMUI 14.2, HarmonyOS 4.2', '6.6" 1256x2760 pixels', '50MP 2160p', '12GB RAM Kirin 9010', '4900mAh 66W 50W')
    ,('HMD', '৳18,000.00', '2024, April 24', 'Android 14', '6.65" 720x1612 pixels', '13MP 1080p', '4/6GB RAM Unisoc T606', '5000mAh Li-Po')
    ,('Vivo', '৳18,000.00', '2024, April 24', 'Android 14, Funtouch 14', '6.72" 1080x2408 pixels', '50MP 2160p', '4-8GB RAM Snapdragon 6 Gen 1', '6000mAh 44W')
    ,('Vivo', '৳25,000.00', '2024, April 27', 'Android 14, OriginOS 4', '6.72" 1080x2408 pixels', '50MP 1080p', '8/12GB RAM Snapdragon 4 Gen 2', '6000mAh 44W')
    ,('Motorola', '৳120,000.00', 'Exp. release 2024, April', 'Android 14', '6.7" 1220x2712 pixels', '50MP 4320p', '12/16GB RAM Snapdragon 8s Gen 3', '5000mAh 125W 50W')
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt 8; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}